# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (F column) and "最低票价" (G column) figures on the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 6966
$ws1.Range("F4").Value = 57
$ws1.Range("F7").Value = 6838
$ws1.Range("F10").Value = 0
$ws1.Range("F12").Value = 108
$ws1.Range("F14").Value = 149
$ws1.Range("F16").Value = 410
$ws1.Range("F18").Value = 39
$ws1.Range("F19").Value = 16
$ws1.Range("G19").Value = 45
$ws1.Range("F20").Value = 5213
$ws1.Range("F22").Value = 164
$ws1.Range("F23").Value = 0
$ws1.Range("F24").Value = 216
$ws1.Range("F25").Value = 232

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 6966
$ws4.Range("F4").Value = 57
$ws4.Range("F5").Value = 454
$ws4.Range("F7").Value = 6838
$ws4.Range("F8").Value = 73
$ws4.Range("F11").Value = 20
$ws4.Range("F14").Value = 149
$ws4.Range("F15").Value = 17
$ws4.Range("F16").Value = 410
$ws4.Range("F18").Value = 39
$ws4.Range("F19").Value = 16
$ws4.Range("G19").Value = 45
$ws4.Range("F21").Value = 5213
$ws4.Range("F22").Value = 0
$ws4.Range("F24").Value = 164
$ws4.Range("F25").Value = 631
$ws4.Range("F26").Value = 216
$ws4.Range("F27").Value = 232
